$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 96.40000000000001
$ws.Range("D2").Value = 13.2
$ws.Range("C3").Value = 39.4
$ws.Range("D3").Value = 99.09999999999999
$ws.Range("C4").Value = 85.5
$ws.Range("D4").Value = 66
$ws.Range("C5").Value = 31.5
$ws.Range("D5").Value = 99.7
